$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: "¿Qué se hizo ayer?" for Paula Andrea Taborda Jaramillo
$ws.Range("C10").Value = "Nada"
$ws.Range("D10").Value = "Asistimos a la reunión y se asignaron tareas en la tabla de casos de uso. Guillermo y yo terminamos el modelo de clases, le fata la cardinalidad"
$ws.Range("E10").Value = "Trabajamos en el diagrama de clases."
$ws.Range("F10").Value = "Nada"
$ws.Range("G10").Value = "Tener en cuenta las correcciones"

# Row 11: "¿Qué se hará hoy?"
$ws.Range("C11").Value = "Reunión para asignar tareas, revisión a lo que hizo Guillermo"
$ws.Range("D11").Value = "Nada"
$ws.Range("E11").Value = "Nada"
$ws.Range("F11").Value = "Mostrar los avances en la clase, anotar las correcciones"
$ws.Range("G11").Value = "Asistir a la reunión"

# Row 12: "¿Qué cosas se oponen?"
$ws.Range("C12").Value = "Ninguna"
$ws.Range("D12").Value = "Semana de parciales"
$ws.Range("E12").Value = "Semana de Parciales"
$ws.Range("F12").Value = "Ninguna"
$ws.Range("G12").Value = "Ninguna"

# Update selection to match the committed state (active cell E12)
$ws.Range("E12").Select()
